# "Test & Gannt fertiggestellt" -- mark the remaining Gantt tasks (rows 36-50)
# as fully completed: copy the start (C) and duration (D) values into the
# "Ist" (actual) start/duration columns E and F, and set the completion
# flag in column G to 1 (100%).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projekt")

$firstRow = 36
$lastRow  = 50

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $start    = $ws.Range("C$r").Value2
    $duration = $ws.Range("D$r").Value2

    $ws.Range("E$r").Value = $start
    $ws.Range("F$r").Value = $duration
    $ws.Range("G$r").Value = 1
}

# Reflect the final cell the author left selected once the plan was completed
$ws.Range("H53").Select()
